# Update gh-pages to output generated at 456a3b4
# Bump the "想去人数" (wanted-to-go count) column F for the events whose
# numbers changed between builds. The "全部类型" sheet aggregates the
# same events (at different row offsets) as the other three sheets, so
# it needs the same updates applied to its own rows.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F18").Value = 563
$ws1.Range("F21").Value = 681
$ws1.Range("F24").Value = 333
$ws1.Range("F27").Value = 697
$ws1.Range("F28").Value = 8029
$ws1.Range("F35").Value = 1670
$ws1.Range("F37").Value = 162

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 68
$ws2.Range("F17").Value = 228

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 769

# 全部类型 (All types) - aggregated view, same events different rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 769
$ws4.Range("F17").Value = 563
$ws4.Range("F21").Value = 68
$ws4.Range("F22").Value = 681
$ws4.Range("F27").Value = 697
$ws4.Range("F28").Value = 8029
